# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new rows 31 and 32) above the existing
# data, pushing the former rows 31-79 down to rows 33-81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 31; this shifts the old
# rows 31..79 down to 33..81 and carries their formatting with them.
$ws.Rows("31:32").Insert()

# ---- New row 31: "Candy White" ----
$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value = 44935
$ws.Cells.Item(31, 5).Value = 15
$ws.Cells.Item(31, 6).Value = "Fruta"
$ws.Cells.Item(31, 7).Value = 100103
$ws.Cells.Item(31, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(31, 9).Value = 100103006
$ws.Cells.Item(31, 10).Value = "Nectarín"
$ws.Cells.Item(31, 11).Value = "Candy White"
$ws.Cells.Item(31, 12).Value = "Segunda"
$ws.Cells.Item(31, 13).Value = 350
$ws.Cells.Item(31, 14).Value = 19000
$ws.Cells.Item(31, 15).Value = 20000
$ws.Cells.Item(31, 16).Value = 19429
$ws.Cells.Item(31, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(31, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(31, 19).Value = 1079
$ws.Cells.Item(31, 20).Value = 18

# ---- New row 32: "Super Queen" ----
$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value = 44935
$ws.Cells.Item(32, 5).Value = 15
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100103
$ws.Cells.Item(32, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(32, 9).Value = 100103006
$ws.Cells.Item(32, 10).Value = "Nectarín"
$ws.Cells.Item(32, 11).Value = "Super Queen"
$ws.Cells.Item(32, 12).Value = "Segunda"
$ws.Cells.Item(32, 13).Value = 500
$ws.Cells.Item(32, 14).Value = 19000
$ws.Cells.Item(32, 15).Value = 20000
$ws.Cells.Item(32, 16).Value = 19400
$ws.Cells.Item(32, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(32, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(32, 19).Value = 1078
$ws.Cells.Item(32, 20).Value = 18
